$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 699
$ws.Range("I2").Value = 465.33334
$ws.Range("K2").Value = 465.33334
$ws.Range("M2").Value = -352.33334
# Row 9
$ws.Range("H9").Value = 112.25
$ws.Range("I9").Value = 68.42856999999999
$ws.Range("J9").Value = 419
$ws.Range("K9").Value = 68.42856999999999
$ws.Range("L9").Value = 419
$ws.Range("M9").Value = 100.57143
$ws.Range("N9").Value = -757
# Row 95
$ws.Range("H95").Value = 80000
$ws.Range("J95").Value = 80000
$ws.Range("L95").Value = 80000
$ws.Range("N95").Value = -85492
# Row 116
$ws.Range("H116").Value = 5288.0625
$ws.Range("I116").Value = 3656.111
$ws.Range("J116").Value = 7386.2856
$ws.Range("K116").Value = 3656.111
$ws.Range("L116").Value = 7386.2856
$ws.Range("M116").Value = -214.1109999999999
$ws.Range("N116").Value = -14270.2856
# Row 128
$ws.Range("H128").Value = 37620.2
$ws.Range("J128").Value = 32974
$ws.Range("L128").Value = 32974
$ws.Range("N128").Value = -42934

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 14500
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15804
# Row 32
$ws.Range("H32").Value = 5670.9
$ws.Range("I32").Value = 3924.913
$ws.Range("K32").Value = 3924.913
$ws.Range("M32").Value = -3637.913
# Row 97
$ws.Range("H97").Value = 1039.1177
$ws.Range("I97").Value = 1074.6
$ws.Range("K97").Value = 1074.6
$ws.Range("M97").Value = -578.5999999999999
# Row 122
$ws.Range("H122").Value = 3177
$ws.Range("I122").Value = 2642
$ws.Range("K122").Value = 7926
$ws.Range("M122").Value = -5476

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4431.075
$ws.Range("I86").Value = 4017.258
$ws.Range("K86").Value = 4017.258
$ws.Range("M86").Value = -2894.258
# Row 89
$ws.Range("H89").Value = 4431.075
$ws.Range("I89").Value = 4017.258
$ws.Range("K89").Value = 20086.29
$ws.Range("M89").Value = -14470.29
# Row 94
$ws.Range("H94").Value = 1081.4375
$ws.Range("I94").Value = 1081.4375
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1081.4375
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -630.4375
$ws.Range("N94").ClearContents()
# Row 105
$ws.Range("H105").Value = 9820.950000000001
$ws.Range("I105").Value = 7941.0312
$ws.Range("K105").Value = 7941.0312
$ws.Range("M105").Value = -6194.0312
# Row 134
$ws.Range("H134").Value = 1851.7142
$ws.Range("I134").Value = 1164.5227
$ws.Range("J134").Value = 4371.4165
$ws.Range("K134").Value = 3493.5681
$ws.Range("L134").Value = 13114.2495
$ws.Range("M134").Value = -958.5681
$ws.Range("N134").Value = -18184.2495

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 23999.666
$ws.Range("J9").Value = 23999.666
$ws.Range("L9").Value = 23999.666
$ws.Range("N9").Value = -24335.666
# Row 31
$ws.Range("H31").Value = 53411.906
$ws.Range("I31").Value = 2874.7778
$ws.Range("J31").Value = 91314.75
$ws.Range("K31").Value = 2874.7778
$ws.Range("L31").Value = 91314.75
$ws.Range("M31").Value = -2579.7778
$ws.Range("N31").Value = -91904.75
# Row 34
$ws.Range("H34").Value = 53411.906
$ws.Range("I34").Value = 2874.7778
$ws.Range("J34").Value = 91314.75
$ws.Range("K34").Value = 2874.7778
$ws.Range("L34").Value = 91314.75
$ws.Range("M34").Value = -2672.7778
$ws.Range("N34").Value = -91718.75
# Row 86
$ws.Range("H86").Value = 4599.9443
$ws.Range("I86").Value = 4174.4546
$ws.Range("J86").Value = 5268.5713
$ws.Range("K86").Value = 4174.4546
$ws.Range("L86").Value = 5268.5713
$ws.Range("M86").Value = -3051.4546
$ws.Range("N86").Value = -7514.5713
# Row 89
$ws.Range("H89").Value = 4599.9443
$ws.Range("I89").Value = 4174.4546
$ws.Range("J89").Value = 5268.5713
$ws.Range("K89").Value = 20872.273
$ws.Range("L89").Value = 26342.8565
$ws.Range("M89").Value = -15256.273
$ws.Range("N89").Value = -37574.85649999999
# Row 94
$ws.Range("H94").Value = 3880.9167
$ws.Range("I94").Value = 3617
$ws.Range("J94").Value = 4069.4285
$ws.Range("K94").Value = 3617
$ws.Range("L94").Value = 4069.4285
$ws.Range("M94").Value = -3166
$ws.Range("N94").Value = -4971.4285
# Row 99
$ws.Range("H99").Value = 1893.5834
$ws.Range("I99").Value = 1679.8182
$ws.Range("K99").Value = 1679.8182
$ws.Range("M99").Value = -181.8181999999999
# Row 107
$ws.Range("H107").Value = 1780.3889
$ws.Range("I107").Value = 1076.1
$ws.Range("K107").Value = 1076.1
$ws.Range("M107").Value = 843.9000000000001
# Row 122
$ws.Range("H122").Value = 5212.2354
$ws.Range("I122").Value = 1601.125
$ws.Range("J122").Value = 8422.111000000001
$ws.Range("K122").Value = 4803.375
$ws.Range("L122").Value = 25266.333
$ws.Range("M122").Value = -2353.375
$ws.Range("N122").Value = -30166.333
# Row 126
$ws.Range("H126").Value = 1893.5834
$ws.Range("I126").Value = 1679.8182
$ws.Range("K126").Value = 5039.4546
$ws.Range("M126").Value = -2569.4546
# Row 133
$ws.Range("H133").Value = 44320.715
$ws.Range("J133").Value = 41408.184
$ws.Range("L133").Value = 41408.184
$ws.Range("N133").Value = -46468.184
# Row 134
$ws.Range("H134").Value = 2242.5264
$ws.Range("I134").Value = 1431.5518
$ws.Range("K134").Value = 4294.6554
$ws.Range("M134").Value = -1759.6554

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 12826759
$ws.Range("J129").Value = 27789616
$ws.Range("L129").Value = 83368848
$ws.Range("N129").Value = -83378848
# Row 131
$ws.Range("H131").Value = 6946206
$ws.Range("I131").Value = 17858158
$ws.Range("K131").Value = 53574474
$ws.Range("M131").Value = -53569434

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 42005
$ws.Range("I20").Value = 42005
$ws.Range("K20").Value = 42005
$ws.Range("M20").Value = -41760
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
# Row 70
$ws.Range("H70").Value = 16599.2
$ws.Range("I70").Value = 14999
$ws.Range("J70").Value = 17666
$ws.Range("K70").Value = 14999
$ws.Range("L70").Value = 17666
$ws.Range("M70").Value = -14729
$ws.Range("N70").Value = -18206
# Row 73
$ws.Range("H73").Value = 16599.2
$ws.Range("I73").Value = 14999
$ws.Range("J73").Value = 17666
$ws.Range("K73").Value = 14999
$ws.Range("L73").Value = 17666
$ws.Range("M73").Value = -14063
$ws.Range("N73").Value = -19538
# Row 92
$ws.Range("H92").Value = 17870.75
$ws.Range("J92").Value = 17870.75
$ws.Range("L92").Value = 17870.75
$ws.Range("N92").Value = -21614.75
# Row 107
$ws.Range("H107").Value = 468.31033
$ws.Range("I107").Value = 462.26086
$ws.Range("J107").Value = 491.5
$ws.Range("K107").Value = 462.26086
$ws.Range("L107").Value = 491.5
$ws.Range("M107").Value = 1457.73914
$ws.Range("N107").Value = -4331.5
# Row 122
$ws.Range("H122").Value = 10012.619
$ws.Range("I122").Value = 10209.4375
$ws.Range("J122").Value = 9382.799999999999
$ws.Range("K122").Value = 30628.3125
$ws.Range("L122").Value = 28148.4
$ws.Range("M122").Value = -28178.3125
$ws.Range("N122").Value = -33048.39999999999
# Row 126
$ws.Range("H126").Value = 3779.6365
$ws.Range("I126").Value = 2135.818
$ws.Range("K126").Value = 6407.454000000001
$ws.Range("M126").Value = -3937.454000000001
# Row 128
$ws.Range("H128").Value = 59998.332
$ws.Range("J128").Value = 59998.332
$ws.Range("L128").Value = 59998.332
$ws.Range("N128").Value = -69958.33199999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 44
$ws.Range("H44").Value = 99097
$ws.Range("J44").Value = 99210.336
$ws.Range("L44").Value = 99210.336
$ws.Range("N44").Value = -100122.336
# Row 69
$ws.Range("H69").Value = 36999
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 36999
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 110
$ws.Range("H110").Value = 39995
$ws.Range("J110").Value = 39995
$ws.Range("L110").Value = 39995
$ws.Range("N110").Value = -48175

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Range("H31").Value = 30000
$ws.Range("J31").Value = 30000
$ws.Range("L31").Value = 30000
$ws.Range("N31").Value = -30696
# Row 132
$ws.Range("H132").Value = 2375.182
$ws.Range("I132").Value = 1726.4231
$ws.Range("J132").Value = 4784.857
$ws.Range("K132").Value = 5179.2693
$ws.Range("L132").Value = 14354.571
$ws.Range("M132").Value = -2649.2693
$ws.Range("N132").Value = -19414.571
# Row 136
$ws.Range("H136").Value = 2991.9412
$ws.Range("I136").Value = 965.45
$ws.Range("J136").Value = 5886.9287
$ws.Range("K136").Value = 2896.35
$ws.Range("L136").Value = 17660.7861
$ws.Range("M136").Value = -346.3500000000004
$ws.Range("N136").Value = -22760.7861
